# Cải tiến thuật toán GA - Xử lý xung đột theo tuần
#
# For every weekly sheet (Tuan_1 .. Tuan_15) the GA re-scheduled two
# classes that were colliding on the same slot:
#   - Row 9  (class CL10, "Ngữ pháp tiếng Anh", Phòng: R105):
#       Slot C2 (15:00-17:00) -> T1 (17:30-19:30)
#       Day   Thứ 2 (col D)   -> Thứ 3 (col E)
#   - Row 10 (class CL05, "Ngữ pháp tiếng Anh"):
#       Slot T1 (17:30-19:30) -> T2 (19:30-21:30)
#       Day   Thứ 4 (col F)   -> Thứ 2 (col D)
#       Room  R102            -> R103
#
# The cell that receives the schedule block must carry the "highlighted"
# style (s=9) while the vacated cell reverts to the plain grid style
# (s=8). We achieve that by copying a neighboring already-plain (s=8)
# cell onto the vacated cell after clearing it.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le 15; $i++) {
    $ws = $wb.Worksheets.Item("Tuan_$i")

    # ---- Row 9: CL10 class moves from Monday(D) to Tuesday(E) slot ----
    $ws.Range("A9").Value = "T1`n(17:30-19:30)"

    # Move the class block (value + highlighted style) from D9 to E9
    $ws.Range("D9").Copy($ws.Range("E9"))

    # Vacate D9: clear its content, then restore the plain (unhighlighted)
    # grid style by copying from a neighboring plain cell.
    $ws.Range("D9").ClearContents()
    $ws.Range("G9").Copy($ws.Range("D9"))

    # ---- Row 10: CL05 class moves from Wednesday(F) to Monday(D) slot,
    #              and its room changes from R102 to R103 ----
    $ws.Range("A10").Value = "T2`n(19:30-21:30)"
    $ws.Range("F10").Value = "Ngữ pháp tiếng Anh`n(Lý thuyết)`nPhòng: R103`nGV: Ngô Văn I"

    # Move the (now room-updated) class block from F10 to D10
    $ws.Range("F10").Copy($ws.Range("D10"))

    # Vacate F10: clear its content, then restore the plain grid style.
    $ws.Range("F10").ClearContents()
    $ws.Range("G10").Copy($ws.Range("F10"))
}
